$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 128
$ws.Range("I61").Value = 135.33333
$ws.Range("J61").Value = 95
$ws.Range("K61").Value = 405.99999
$ws.Range("L61").Value = 285
$ws.Range("M61").Value = -233.99999
$ws.Range("N61").Value = -629

$ws.Range("H111").Value = 2855.5
$ws.Range("I111").Value = 2516.3333
$ws.Range("J111").Value = 3194.6667
$ws.Range("K111").Value = 7548.999899999999
$ws.Range("L111").Value = 9584.000100000001
$ws.Range("M111").Value = -4481.999899999999
$ws.Range("N111").Value = -15718.0001

$ws.Range("H121").Value = 3639
$ws.Range("J121").Value = 7605
$ws.Range("L121").Value = 22815
$ws.Range("N121").Value = -26309

$ws.Range("H129").Value = 1199.7297
$ws.Range("I129").Value = 782.25
$ws.Range("J129").Value = 1250.3334
$ws.Range("K129").Value = 2346.75
$ws.Range("L129").Value = 3751.0002
$ws.Range("M129").Value = 2653.25
$ws.Range("N129").Value = -13751.0002

$ws.Range("H138").Value = 1953.1875
$ws.Range("I138").Value = 1693.8846
$ws.Range("J138").Value = 2259.6365
$ws.Range("K138").Value = 5081.6538
$ws.Range("L138").Value = 6778.9095
$ws.Range("M138").Value = 58.34619999999995
$ws.Range("N138").Value = -17058.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 349.66666
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 274.5
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 274.5
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = -506.5

$ws.Range("H5").Value = 83.42856999999999
$ws.Range("I5").Value = 80.666664
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 80.666664
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 31.333336
$ws.Range("N5").Value = -324

$ws.Range("H45").Value = 1138.5294
$ws.Range("I45").Value = 972
$ws.Range("K45").Value = 972
$ws.Range("M45").Value = -595

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 83.42856999999999
$ws.Range("I4").Value = 80.666664
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 80.666664
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 34.333336
$ws.Range("N4").Value = -330

$ws.Range("H20").Value = 1591.6666
$ws.Range("I20").Value = 1004.6429
$ws.Range("J20").Value = 2223.8462
$ws.Range("K20").Value = 1004.6429
$ws.Range("L20").Value = 2223.8462
$ws.Range("M20").Value = -757.6429000000001
$ws.Range("N20").Value = -2717.8462

$ws.Range("H80").Value = 294.96155
$ws.Range("I80").Value = 266.69232
$ws.Range("J80").Value = 323.23077
$ws.Range("K80").Value = 266.69232
$ws.Range("L80").Value = 323.23077
$ws.Range("M80").Value = 731.30768
$ws.Range("N80").Value = -2319.23077

$ws.Range("H83").Value = 294.96155
$ws.Range("I83").Value = 266.69232
$ws.Range("J83").Value = 323.23077
$ws.Range("K83").Value = 1333.4616
$ws.Range("L83").Value = 1616.15385
$ws.Range("M83").Value = 3658.5384
$ws.Range("N83").Value = -11600.15385

$ws.Range("H94").Value = 2779827.8
$ws.Range("I94").Value = 5557588.5
$ws.Range("J94").Value = 2066.6667
$ws.Range("K94").Value = 5557588.5
$ws.Range("L94").Value = 2066.6667
$ws.Range("M94").Value = -5557137.5
$ws.Range("N94").Value = -2968.6667

$ws.Range("H99").Value = 1971.1111
$ws.Range("I99").Value = 1225
$ws.Range("J99").Value = 2568
$ws.Range("K99").Value = 1225
$ws.Range("L99").Value = 2568
$ws.Range("M99").Value = 273
$ws.Range("N99").Value = -5564

$ws.Range("H107").Value = 1100.2
$ws.Range("I107").Value = 1100.2
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1100.2
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 819.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7354.2856
$ws.Range("J4").Value = 8413.333000000001
$ws.Range("L4").Value = 8413.333000000001
$ws.Range("N4").Value = -8637.333000000001

$ws.Range("H63").Value = 72271
$ws.Range("J63").Value = 72271
$ws.Range("L63").Value = 72271
$ws.Range("N63").Value = -73643

$ws.Range("H66").Value = 72271
$ws.Range("J66").Value = 72271
$ws.Range("L66").Value = 216813
$ws.Range("N66").Value = -223677

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 8623.333000000001
$ws.Range("J105").Value = 8950.362999999999
$ws.Range("L105").Value = 26851.089
$ws.Range("N105").Value = -32093.089

$ws.Range("H113").Value = 506.92307
$ws.Range("I113").Value = 421.25
$ws.Range("J113").Value = 644
$ws.Range("K113").Value = 1263.75
$ws.Range("L113").Value = 1932
$ws.Range("M113").Value = 906.25
$ws.Range("N113").Value = -6272

$ws.Range("H129").Value = 2771.9333
$ws.Range("I129").Value = 1097.7778
$ws.Range("J129").Value = 5283.1665
$ws.Range("K129").Value = 3293.3334
$ws.Range("L129").Value = 15849.4995
$ws.Range("M129").Value = 1706.6666
$ws.Range("N129").Value = -25849.4995

$ws.Range("H131").Value = 916.0513
$ws.Range("J131").Value = 1214.64
$ws.Range("L131").Value = 3643.92
$ws.Range("N131").Value = -13723.92

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 53.166668
$ws.Range("I2").Value = 53
$ws.Range("J2").Value = 53.25
$ws.Range("K2").Value = 53
$ws.Range("L2").Value = 53.25
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = -279.25

$ws.Range("H95").Value = 9508.4
$ws.Range("J95").Value = 9508.4
$ws.Range("L95").Value = 9508.4
$ws.Range("N95").Value = -15000.4

$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 37590.5
$ws.Range("I132").Value = 51233.906
$ws.Range("J132").Value = 5755.8887
$ws.Range("K132").Value = 153701.718
$ws.Range("L132").Value = 17267.6661
$ws.Range("M132").Value = -151171.718
$ws.Range("N132").Value = -22327.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1385
$ws.Range("I132").Value = 1156.9796
$ws.Range("J132").Value = 1892.8636
$ws.Range("K132").Value = 3470.9388
$ws.Range("L132").Value = 5678.5908
$ws.Range("M132").Value = -940.9387999999999
$ws.Range("N132").Value = -10738.5908
